$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update example clinical info: neoadjuvant/adjuvant/second-line/third-line regimens
$ws.Range("B2").Value = "Irinotecan"
$ws.Range("E2").Value = "Capecitabine"
$ws.Range("K2").Value = "5FU, Leucovorin"
$ws.Range("N2").Value = "Oxaliplatin, Other"

# Move active selection to reflect the edited cell (N2), matching the saved view state
$ws.Range("N2").Select()
